$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate stats after trade #88 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.68    # Current Capital
$summary.Range("B4").Value = -0.33      # Total P&L $
$summary.Range("B5").Value = -0.07      # Total P&L %
$summary.Range("B6").Value = 88         # Total Trades
$summary.Range("B8").Value = 47         # Losing Trades
$summary.Range("B9").Value = 32.95      # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking row (row 4) stats.
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.68       # Capital
$status.Range("D4").Value = 88          # Trades
$status.Range("E4").Value = -0.33       # P&L $
$status.Range("F4").Value = -0.32       # P&L %
$status.Range("G4").Value = 32.95       # Win Rate %

# ---------------------------------------------------------------------------
# New closed trade (#88) appended as row 89 on both the "All Trades" and
# "MarketMaking" logs (identical content in this workbook).
# ---------------------------------------------------------------------------
function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value = 88

    # Columns B/C hold plain-text date/time strings in this workbook. Column
    # B ("2026-02-17") looks like a date to Excel's auto-detection, so force
    # text formatting for that one cell, write the literal value, then put
    # the cell style back to Normal so no stray number format lingers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "15:53:41"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.881083
    $ws.Cells.Item($row, 7).Value = 0.64
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -27.3622
    $ws.Cells.Item($row, 10).Value = -0.24
    $ws.Cells.Item($row, 11).Value = 99.68
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 89

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 89

Write-Host "Trade #88 recorded; Summary and Strategy Status refreshed."
